$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.528.57'
$ws.Range('E2').Value = '  -1.80%  '

# Row 3
$ws.Range('D3').Value = '1.921.46'
$ws.Range('E3').Value = '  -2.04%  '

# Row 4
$ws.Range('D4').Value = "'" + '1.005'
$ws.Range('E4').Value = '  +0.39%  '

# Row 5
$ws.Range('D5').Value = "'" + '239.10'
$ws.Range('E5').Value = '  -2.76%  '

# Row 6
$ws.Range('D6').Value = "'" + '1.004'
$ws.Range('E6').Value = '  +0.36%  '

# Row 7
$ws.Range('D7').Value = "'" + '0.4796'
$ws.Range('E7').Value = '  -1.84%  '

# Row 8
$ws.Range('D8').Value = "'" + '0.2878'
$ws.Range('E8').Value = '  -3.00%  '

# Row 9
$ws.Range('D9').Value = "'" + '0.06704'
$ws.Range('E9').Value = '  -1.88%  '

# Row 10
$ws.Range('D10').Value = "'" + '18.81'
$ws.Range('E10').Value = '  -2.26%  '

# Row 11
$ws.Range('D11').Value = "'" + '104.24'
$ws.Range('E11').Value = '  -2.89%  '

# Row 12
$ws.Range('D12').Value = '1.924.49'
$ws.Range('E12').Value = '  -2.64%  '

# Row 13
$ws.Range('D13').Value = "'" + '0.07752'
$ws.Range('E13').Value = '  -1.05%  '

# Row 14
$ws.Range('D14').Value = "'" + '5.245'
$ws.Range('E14').Value = '  -4.66%  '

# Row 15
$ws.Range('D15').Value = "'" + '0.6821'
$ws.Range('E15').Value = '  -3.41%  '

# Row 16
$ws.Range('D16').Value = "'" + '265.58'
$ws.Range('E16').Value = '  -6.93%  '

# Row 17
$ws.Range('D17').Value = '30.573.29'
$ws.Range('E17').Value = '  -1.68%  '

# Row 18
$ws.Range('D18').Value = "'" + '1.004'

# Row 19
$ws.Range('D19').Value = "'" + '0.000007527'
$ws.Range('E19').Value = '  -2.39%  '

# Row 20
$ws.Range('D20').Value = "'" + '12.73'
$ws.Range('E20').Value = '  -3.92%  '

# Row 21
$ws.Range('D21').Value = "'" + '5.435'
$ws.Range('E21').Value = '  -1.69%  '

# Row 22
$ws.Range('D22').Value = "'" + '1.005'
$ws.Range('E22').Value = '  +0.46%  '

# Row 23
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'" + '6.353'
$ws.Range('E23').Value = '  -2.39%  '

# Row 24
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = "'" + '9.646'
$ws.Range('E24').Value = '  -1.86%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'" + '163.66'
$ws.Range('E25').Value = '  -3.35%  '

# Row 26
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = "'" + '19.03'
$ws.Range('E26').Value = '  -5.14%  '

# Row 27
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = "'" + '2.096'
$ws.Range('E27').Value = '  -5.26%  '

# Row 28
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = "'" + '0.1022'
$ws.Range('E28').Value = '  -3.39%  '

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'" + '1.389'
$ws.Range('E29').Value = '  -1.48%  '

# Row 30
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'" + '4.603'
$ws.Range('E30').Value = '  -0.42%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'" + '1.519'
$ws.Range('E31').Value = '  -4.36%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'" + '4.264'
$ws.Range('E32').Value = '  -4.19%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'" + '0.04757'
$ws.Range('E33').Value = '  -3.73%  '

# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'" + '0.7394'
$ws.Range('E34').Value = '  -3.22%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'" + '1.120'
$ws.Range('E35').Value = '  -4.80%  '

# Row 36
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').Value = "'" + '1.003'
$ws.Range('E36').Value = '  +0.28%  '

# Row 37
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = "'" + '2.680'
$ws.Range('E37').Value = '  -1.97%  '

# Row 38
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'" + '0.01940'
$ws.Range('E38').Value = '  -3.91%  '

# Row 39
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'" + '2.643'
$ws.Range('E39').Value = '  -2.22%  '

# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'" + '6.341'
$ws.Range('E40').Value = '  -3.23%  '

# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = "'" + '75.57'
$ws.Range('E41').Value = '  -2.90%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = "'" + '2.001'
$ws.Range('E42').Value = '  -5.64%  '

# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'" + '0.8624'
$ws.Range('E43').Value = '  -3.49%  '

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'" + '106.19'
$ws.Range('E44').Value = '  -2.65%  '

# Row 45
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = "'" + '0.4288'
$ws.Range('E45').Value = '  -4.52%  '

# Row 46
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = "'" + '1.004'
$ws.Range('E46').Value = '  +0.25%  '

# Row 47
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = "'" + '7.563'
$ws.Range('E47').Value = '  -7.70%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = "'" + '997.82'
$ws.Range('E48').Value = '  -0.98%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'" + '0.1203'
$ws.Range('E49').Value = '  -4.81%  '

# Row 50
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'" + '35.22'
$ws.Range('E50').Value = '  -2.23%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'" + '9.007'
$ws.Range('E51').Value = '  -4.07%  '
